# Generate Report for Handoff
#
# The "e654846b-93f6-4791-8f6c-75839332b73c.md" file has moved from
# "In Translation" to "Ready for handoff" in both locales, so update the
# Overview sheet as well as the per-locale (zh-cn / de-de) detail sheets
# with the new status and handoff timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for e654846b-...md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-17-11 09:17:33"

# --- zh-cn sheet: row for e654846b-...md (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-11 09:17:30"

# --- de-de sheet: row for e654846b-...md (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-11 09:17:33"
